$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CRMAccuracyData")

$ws.Range("A84").Value = 20220504
$ws.Range("B84").Value = 2200.1692800000001
$ws.Range("C84").Value = 2224.4699999999998
$ws.Range("D84").Formula = "=100*(B84-C84)/C84"
$ws.Range("E84").Value = 180
$ws.Range("F84").Value = "CRM OPENED 20220427"
